# Update "想去人数" (interested-people count) figures in column F
# for the "展览" and "全部类型" worksheets, matching the refreshed
# scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibition) ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 56
$wsExpo.Range("F4").Value  = 76
$wsExpo.Range("F5").Value  = 35
$wsExpo.Range("F6").Value  = 35
$wsExpo.Range("F7").Value  = 2637
$wsExpo.Range("F8").Value  = 1150
$wsExpo.Range("F9").Value  = 234
$wsExpo.Range("F11").Value = 5931
$wsExpo.Range("F13").Value = 233
$wsExpo.Range("F14").Value = 582
$wsExpo.Range("F15").Value = 11597
$wsExpo.Range("F16").Value = 11818

# ---- Sheet "全部类型" (All types) ----
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 56
$wsAll.Range("F4").Value  = 76
$wsAll.Range("F5").Value  = 35
$wsAll.Range("F6").Value  = 35
$wsAll.Range("F7").Value  = 2637
$wsAll.Range("F9").Value  = 1150
$wsAll.Range("F10").Value = 234
$wsAll.Range("F12").Value = 5931
$wsAll.Range("F14").Value = 233
$wsAll.Range("F15").Value = 582
$wsAll.Range("F16").Value = 11597
$wsAll.Range("F17").Value = 11818
